$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: RandomForestRegressor (name unchanged)
$ws.Range("B3").Value = 0.9969767486873989
$ws.Range("C3").Value = 0.997255990073923
$ws.Range("D3").Value = 0.996864189913787

# Row 4: GradientBoostingRegressor -> DecisionTreeRegressor
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.9968166501530092
$ws.Range("C4").Value = 0.9971960380038718
$ws.Range("D4").Value = 0.9971574624206325

# Row 5: AdaBoostRegressor -> MLPRegressor
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.9982917352888737
$ws.Range("C5").Value = 0.9980344099994461
$ws.Range("D5").Value = 0.9979167009397355
